# Update sample data on the "All" sheet:
#   First Name : Loren -> Charles
#   Last Name  : Velasquez -> Tiberius
#   Job        : Software Engineer -> Data Engineer

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All")

$ws.Range("B1").Value = "Charles"
$ws.Range("B4").Value = "Data Engineer"
$ws.Range("B2").Value = "Tiberius"

# Reflect the new active cell / selection recorded for this sheet.
$ws.Activate()
$ws.Range("D13").Select()
